# Updated cryptos list on Tue Nov 19 21:43:48 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "92.143.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.086.94"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.35%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.23"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "608.90"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.08"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.94%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.387"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.081.65"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.726"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.44%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.332.71"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.92"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.41"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.675.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.105.69"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.74"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.52"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.70"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.28"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "441.06"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000193"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.66"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "85.62"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.56"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.253.56"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.96%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.130"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.168"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.225"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.05"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.83%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.79"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.156"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -8.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.74"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.14%  "

$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.88"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.75%  "

$ws.Range("B39").Value = "MantraDAO"
$ws.Range("C39").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.85"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "479.04"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.88"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.58%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.428"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.67%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.27"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "164.64"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.27%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.680"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.36"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.22%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.41%  "
